$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B19 was stored as an inline string "5"; convert to a real number
$ws.Range("B19").Value = 5

# Add new row 20 with data
$ws.Range("A20").Value = "Ying Tang"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "3"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "无"
$ws.Range("D20").Value = "DIS"
$ws.Range("E20").Value = "MET"
$ws.Range("F20").Value = "2bb8b329-99fa-4c06-a5b4-7897e3cce401"
$ws.Range("G20").Value = "S1PWi_lC-_annotated.xlsx"
$ws.Range("H20").Value = "Each network is trained with 50 epochs."
